$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28
$ws.Range("A28").Value = 'Mar'
$ws.Range("B28").Value = 'Jan'
$ws.Range("C28").Value = 'Apr'
$ws.Range("D28").Value = 'Sep'
$ws.Range("E28").Value = 'Tampa, Florida'
$ws.Range("F28").Value = 'Caribbean'
$ws.Range("G28").Value = '6 - 8'
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = '6'
$ws.Range("H28").Style = "Normal"

# Row 29
$ws.Range("A29").Value = 'Mar'
$ws.Range("B29").Value = 'Jan'
$ws.Range("C29").Value = 'Apr'
$ws.Range("D29").Value = 'Sep'
$ws.Range("E29").Value = 'Tampa, Florida'
$ws.Range("F29").Value = 'Caribbean'
$ws.Range("G29").Value = '6 - 8'
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = '6'
$ws.Range("H29").Style = "Normal"
$ws.Range("I29").Value = 'Sep'
$ws.Range("J29").Value = 'Dec'
$ws.Range("K29").Value = 'Mar'
$ws.Range("L29").Value = 'May'
$ws.Range("M29").Value = 'San Juan, Puerto Rico'
$ws.Range("N29").Value = 'Caribbean'
$ws.Range("O29").Value = '2 - 5'
$ws.Range("P29").NumberFormat = "@"
$ws.Range("P29").Value = '2'
$ws.Range("P29").Style = "Normal"

# Row 30
$ws.Range("A30").Value = 'Nov'
$ws.Range("B30").Value = 'Feb'
$ws.Range("C30").Value = 'Apr'
$ws.Range("D30").Value = 'Aug'
$ws.Range("E30").Value = 'San Juan, Puerto Rico'
$ws.Range("F30").Value = 'Caribbean'
$ws.Range("G30").Value = '6 - 8'
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = '8'
$ws.Range("H30").Style = "Normal"

# Row 31
$ws.Range("A31").Value = 'Nov'
$ws.Range("B31").Value = 'Feb'
$ws.Range("C31").Value = 'Apr'
$ws.Range("D31").Value = 'Aug'
$ws.Range("E31").Value = 'San Juan, Puerto Rico'
$ws.Range("F31").Value = 'Caribbean'
$ws.Range("G31").Value = '6 - 8'
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = '8'
$ws.Range("H31").Style = "Normal"
$ws.Range("I31").Value = 'Dec'
$ws.Range("J31").Value = 'Aug'
$ws.Range("K31").Value = 'Feb'
$ws.Range("L31").Value = 'Mar'
$ws.Range("M31").Value = 'San Juan, Puerto Rico'
$ws.Range("N31").Value = 'Caribbean'
$ws.Range("O31").Value = '6 - 8'
$ws.Range("P31").NumberFormat = "@"
$ws.Range("P31").Value = '10'
$ws.Range("P31").Style = "Normal"
